$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: I1 "I0", J1 "IF" ---
# Copy H1's formatting (bold font, border, center/top alignment) onto I1:J1
# so the new header cells reuse the existing style index (s="1") instead of
# creating a new duplicate style entry.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows 2-69: new columns I (I0) and J (IF) ---
$iVals = @(7, 8, 9, 9, 8, 10, 9, 9, 8, 9, 9, 9, 8, 10, 9, 9, 9, 11, 9, 9, 9, 10, 10, 9, 9, 9, 9, 9, 9, 9, 8, 9, 9, 9, 9, 9, 9, 9, 9, 9, 8, 9, 9, 9, 8, 9, 9, 9, 9, 9, 9, 8, 9, 9, 10, 9, 9, 9, 9, 9, 9, 9, 9, 9, 5, 5, 4, 2)
$jVals = @(8, 8, 9, 9, 9, 11, 9, 9, 9, 10, 9, 9, 9, 10, 9, 9, 9, 11, 9, 9, 9, 10, 10, 9, 9, 9, 9, 9, 9, 9, 9, 9, 9, 9, 9, 9, 9, 9, 9, 9, 9, 10, 9, 9, 8, 9, 9, 9, 9, 9, 9, 9, 9, 9, 10, 9, 9, 9, 9, 9, 9, 9, 9, 9, 5, 5, 4, 2)

for ($k = 0; $k -lt $iVals.Length; $k++) {
    $r = $k + 2
    $ws.Cells.Item($r, 9).Value = $iVals[$k]
    $ws.Cells.Item($r, 10).Value = $jVals[$k]
}
